# Weekly update: insert a new price record row for
# "Feria Lagunitas de Puerto Montt - Haba" just after the current top
# (most recent) data row, shifting the existing historical rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 129 (pushes rows 129..158 down to 130..159)
$ws.Rows.Item(129).Insert()

# Populate the new row 129 with the latest weekly record
$ws.Cells.Item(129, 1).Value  = 4
$ws.Cells.Item(129, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(129, 3).Value  = "Los Lagos"
$ws.Cells.Item(129, 4).Value  = 45173
$ws.Cells.Item(129, 5).Value  = 10
$ws.Cells.Item(129, 6).Value  = 100112026
$ws.Cells.Item(129, 7).Value  = "Haba"
$ws.Cells.Item(129, 8).Value  = "Sin especificar"
$ws.Cells.Item(129, 9).Value  = "Primera"
$ws.Cells.Item(129, 10).Value = 40
$ws.Cells.Item(129, 11).Value = 16000
$ws.Cells.Item(129, 12).Value = 16000
$ws.Cells.Item(129, 13).Value = 16000
$ws.Cells.Item(129, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(129, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(129, 16).Value = 640
$ws.Cells.Item(129, 17).Value = 25
$ws.Cells.Item(129, 18).Value = "Hortaliza"
